# BOT; UPDATE DATA
# Applies the daily data-refresh edit described by the commit diff:
#  - "all" sheet:   insert a new data row (43957 / 2020-05-06) before the
#                     trailing total row, push the total row down one row.
#  - "kobe" sheet:   correct H83 (8 -> 9) and populate the previously-blank
#                     row 84 with the new day's figures.
#  - "other" sheet:  populate the previously-blank row 59 with the new
#                     day's figures (copying row 58's number formats).
#  - shared footnote string: extend the list of out-of-city patient
#                     numbers and drop the now-stale "all reported" clause.

$wb = $excel.ActiveWorkbook

$newFootnote = "※　24・34・53・58・59・60・158・161・163・192・237・248・268・272例目（計14件）は市外在住者です。"

# ---------------------------------------------------------------------
# Sheet "all": insert new row 29 (2020-05-06 data), shifting the old
# trailing "total" row from 29 -> 30.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
[void]$wsAll.Rows.Item(29).Insert()

$wsAll.Range("A29").Value = 43957
$wsAll.Range("B29").Value = 272
$wsAll.Range("C29").Value = 268
$wsAll.Range("D29").Value = 102
$wsAll.Range("E29").Value = 92
$wsAll.Range("F29").Value = 10
$wsAll.Range("G29").Value = 7
$wsAll.Range("H29").Value = 159

# Row 30 now holds the shared "total" footnote cell; refresh its text.
$wsAll.Range("B30").Value = $newFootnote

$wsAll.Activate()
[void]$wsAll.Range("B30").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 2

# ---------------------------------------------------------------------
# Sheet "kobe": fix H83, fill in row 84 with the new day's data.
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("H83").Value = 9

$wsKobe.Range("A84").Value = 43957
$wsKobe.Range("B84").Value = 72
$wsKobe.Range("C84").Value = 2351
$wsKobe.Range("D84").Value = 2
$wsKobe.Range("E84").Value = 272
$wsKobe.Range("F84").Value = 97
$wsKobe.Range("G84").Value = 88
$wsKobe.Range("H84").Value = 9
$wsKobe.Range("I84").Value = 7
$wsKobe.Range("J84").Value = 152

# Row 85 holds the shared "total" footnote cell; refresh its text too.
$wsKobe.Range("B85").Value = $newFootnote

$wsKobe.Activate()
[void]$wsKobe.Range("B86").Select()
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 2

# ---------------------------------------------------------------------
# Sheet "other": fill in row 59 with the new day's data, matching row
# 58's number formats (the row was previously blank/unformatted data).
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

[void]$wsOther.Range("A58:I58").Copy()
[void]$wsOther.Range("A59:I59").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsOther.Range("A59").Value = 43957
$wsOther.Range("B59").Value = 0
$wsOther.Range("C59").Value = 12
$wsOther.Range("D59").Value = 5
$wsOther.Range("E59").Value = 4
$wsOther.Range("F59").Value = 1
$wsOther.Range("G59").Value = 0
$wsOther.Range("H59").Value = 7

$wsOther.Activate()
[void]$wsOther.Range("J59").Select()
$excel.ActiveWindow.ScrollRow = 56
$excel.ActiveWindow.ScrollColumn = 5

# Re-select the originally active sheet/cell.
$wsAll.Activate()
